$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '94.194.28'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.92%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.423.93'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.18%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.48'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -6.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '640.82'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.44'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.403'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -5.53%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.964'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -7.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.420.26'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.11%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -5.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.34'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.16'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.207.43'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.51%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.068.05'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.38%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000250'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.28'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -6.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.432.85'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.44'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.41'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.497'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -5.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '496.44'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.21'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000192'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -4.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.30'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -11.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '91.21'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -6.03%  '
$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.98'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.65%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.612.19'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '11.63'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.94%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.74'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +6.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.137'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '29.64'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.550'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '545.32'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.70%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.62'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.43'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -6.21%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.150'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.916'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +6.12%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.71'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.50%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.33'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.14%  '
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.54'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0408'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.84%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.55'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.95%  '
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '54.28'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.91%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.16'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -7.98%  '
